# Updated cryptos list (price / 1h volume change refresh, plus a few row re-orderings).
#
# Price values (column D) are re-entered with a leading apostrophe so Excel keeps them
# as literal text (matching the source data, which mixes "1,234.56"-style thousands
# separators, trailing zeros, leading zeros, and one subscript-notation price with plain
# decimals) instead of silently re-interpreting/normalizing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "'66.881.22"
$ws.Range("E2").Value2 = "  +3.42%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "'3.842.27"
$ws.Range("E3").Value2 = "  +5.08%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value2 = "  -0.28%  "

# Row 5: BNB
$ws.Range("D5").Value = "'422.46"
$ws.Range("E5").Value2 = "  +3.74%  "

# Row 6: Solana
$ws.Range("D6").Value = "'128.88"
$ws.Range("E6").Value2 = "  -3.58%  "

# Row 7: LidoStakedEther
$ws.Range("D7").Value = "'3.834.19"
$ws.Range("E7").Value2 = "  +4.94%  "

# Row 8: XRP
$ws.Range("E8").Value2 = "  -2.32%  "

# Row 9: USDC
$ws.Range("D9").Value = "'0.998"
$ws.Range("E9").Value2 = "  -0.23%  "

# Row 10: Cardano
$ws.Range("E10").Value2 = "  -1.24%  "

# Row 11: Dogecoin
$ws.Range("E11").Value2 = "  -4.64%  "

# Row 12: ShibaInu
$ws.Range("D12").Value = "'0.0000333"
$ws.Range("E12").Value2 = "  +0.55%  "

# Row 13: Avalanche
$ws.Range("D13").Value = "'40.73"
$ws.Range("E13").Value2 = "  -3.03%  "

# Row 14: Polkadot
$ws.Range("D14").Value = "'10.35"
$ws.Range("E14").Value2 = "  +3.65%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "'4.444.07"
$ws.Range("E15").Value2 = "  +4.75%  "

# Row 16: Uniswap
$ws.Range("D16").Value = "'15.61"
$ws.Range("E16").Value2 = "  +14.60%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "'3.857.55"
$ws.Range("E17").Value2 = "  +5.38%  "

# Row 18: TRON
$ws.Range("E18").Value2 = "  -0.64%  "

# Row 19: Chainlink
$ws.Range("D19").Value = "'19.83"
$ws.Range("E19").Value2 = "  -1.08%  "

# Row 20: WrappedBTC
$ws.Range("D20").Value = "'67.088.20"
$ws.Range("E20").Value2 = "  +3.43%  "

# Row 21: Polygon
$ws.Range("E21").Value2 = "  -0.74%  "

# Row 22: BitcoinCash
$ws.Range("D22").Value = "'408.52"
$ws.Range("E22").Value2 = "  -3.21%  "

# Row 23: InternetComputer(DFINITY)
$ws.Range("D23").Value = "'14.89"
$ws.Range("E23").Value2 = "  -1.96%  "

# Row 24: Litecoin
$ws.Range("D24").Value = "'84.07"
$ws.Range("E24").Value2 = "  -2.20%  "

# Row 25: ImmutableX
$ws.Range("D25").Value = "'3.03"

# Row 26: EthereumClassic
$ws.Range("D26").Value = "'37.37"
$ws.Range("E26").Value2 = "  +4.27%  "

# Row 27: Filecoin
$ws.Range("E27").Value2 = "  +5.98%  "

# Row 28: PancakeSwap
$ws.Range("D28").Value = "'3.24"
$ws.Range("E28").Value2 = "  +0.67%  "

# Row 29: RenderToken
$ws.Range("D29").Value = "'9.49"
$ws.Range("E29").Value2 = "  +36.17%  "

# Row 30: LEO
$ws.Range("D30").Value = "'5.40"
$ws.Range("E30").Value2 = "  +5.17%  "

# Row 31: Bittensor
$ws.Range("D31").Value = "'745.57"
$ws.Range("E31").Value2 = "  +9.15%  "

# Row 32: Cosmos
$ws.Range("D32").Value = "'13.13"
$ws.Range("E32").Value2 = "  +2.63%  "

# Row 33: Hedera
$ws.Range("D33").Value = "'0.122"
$ws.Range("E33").Value2 = "  +3.64%  "

# Row 34: Toncoin
$ws.Range("D34").Value = "'2.77"
$ws.Range("E34").Value2 = "  +1.81%  "

# Row 35: Dai
$ws.Range("E35").Value2 = "  -0.15%  "

# Row 36: Kaspa
$ws.Range("E36").Value2 = "  -6.07%  "

# Row 37: InjectiveProtocol
$ws.Range("D37").Value = "'38.34"
$ws.Range("E37").Value2 = "  -7.68%  "

# Row 38: OKB (was NEARProtocol; rows 38-41 reordered)
$ws.Range("B38").Value2 = "OKB"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "'55.48"
$ws.Range("E38").Value2 = "  -0.76%  "

# Row 39: NEARProtocol (was OKB)
$ws.Range("B39").Value2 = "NEARProtocol"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'5.47"
$ws.Range("E39").Value2 = "  +23.00%  "

# Row 40: VeChain (was PEPE)
$ws.Range("B40").Value2 = "VeChain"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0456"
$ws.Range("E40").Value2 = "  -2.37%  "

# Row 41: PEPE (was VeChain)
$ws.Range("B41").Value2 = "PEPE"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "'0.0₃0724"
$ws.Range("E41").Value2 = "  +9.16%  "

# Row 42: ThetaToken
$ws.Range("D42").Value = "'2.88"
$ws.Range("E42").Value2 = "  -2.23%  "

# Row 43: FirstDigitalUSD
$ws.Range("E43").Value2 = "  +0.51%  "

# Row 44: LidoDAOToken
$ws.Range("E44").Value2 = "  +1.14%  "

# Row 45: Stellar
$ws.Range("E45").Value2 = "  -4.79%  "

# Row 46: TheGraph
$ws.Range("E46").Value2 = "  +8.55%  "

# Row 47: ApeXProtocol
$ws.Range("D47").Value = "'3.13"
$ws.Range("E47").Value2 = "  -0.27%  "

# Row 48: ARBITRUM
$ws.Range("E48").Value2 = "  -2.13%  "

# Row 49: Monero
$ws.Range("D49").Value = "'140.74"
$ws.Range("E49").Value2 = "  -2.52%  "

# Row 50: Stacks
$ws.Range("D50").Value = "'2.80"
$ws.Range("E50").Value2 = "  -0.91%  "

# Row 51: WEMIXToken
$ws.Range("E51").Value2 = "  +0.63%  "
